$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 927.2174
$ws.Range("J17").Value = 862.0454999999999
$ws.Range("L17").Value = 2586.1365
$ws.Range("N17").Value = -2922.1365

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1540.875
$ws.Range("I92").Value = 1591.5834
$ws.Range("K92").Value = 1591.5834
$ws.Range("M92").Value = -343.5834

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 563.5
$ws.Range("I96").Value = 563.5
$ws.Range("K96").Value = 1690.5
$ws.Range("M96").Value = -317.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 6228.143
$ws.Range("I106").Value = 5938
$ws.Range("K106").Value = 5938
$ws.Range("M106").Value = -5307

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 12916.667
$ws.Range("I113").Value = 9375
$ws.Range("K113").Value = 9375
$ws.Range("M113").Value = -6121

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 16928.54
$ws.Range("I125").Value = 26496.375
$ws.Range("J125").Value = 1620
$ws.Range("K125").Value = 238467.375
$ws.Range("L125").Value = 14580
$ws.Range("M125").Value = -236007.375
$ws.Range("N125").Value = -19500

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 14409.571
$ws.Range("I132").Value = 2100.0857
$ws.Range("K132").Value = 6300.257100000001
$ws.Range("M132").Value = -3770.257100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5166.537
$ws.Range("I32").Value = 5631.8296
$ws.Range("J32").Value = 2042.4286
$ws.Range("K32").Value = 5631.8296
$ws.Range("L32").Value = 2042.4286
$ws.Range("M32").Value = -5344.8296
$ws.Range("N32").Value = -2616.4286

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 79186.5
$ws.Range("J55").Value = 111098.4
$ws.Range("L55").Value = 111098.4
$ws.Range("N55").Value = -111728.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 10548.223
$ws.Range("I61").Value = 11933.429
$ws.Range("J61").Value = 5700
$ws.Range("K61").Value = 11933.429
$ws.Range("L61").Value = 5700
$ws.Range("M61").Value = -11721.429
$ws.Range("N61").Value = -6124

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 3389.4546
$ws.Range("I97").Value = 1537.1111
$ws.Range("J97").Value = 11725
$ws.Range("K97").Value = 1537.1111
$ws.Range("L97").Value = 11725
$ws.Range("M97").Value = -1041.1111
$ws.Range("N97").Value = -12717

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 10548.223
$ws.Range("I136").Value = 11933.429
$ws.Range("K136").Value = 35800.287
$ws.Range("L136").Value = 17100
$ws.Range("M136").Value = -33250.287
$ws.Range("N136").Value = -22200

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 300
$ws.Range("J22").Value = 250
$ws.Range("L22").Value = 250
$ws.Range("N22").Value = -596

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 35000
$ws.Range("I54").Value = 35000
$ws.Range("K54").Value = 35000
$ws.Range("M54").Value = -34516

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2194
$ws.Range("I94").Value = 868.7059
$ws.Range("J94").Value = 6700
$ws.Range("K94").Value = 868.7059
$ws.Range("L94").Value = 6700
$ws.Range("M94").Value = -417.7059
$ws.Range("N94").Value = -7602

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 6767.407
$ws.Range("I107").Value = 6590.1055
$ws.Range("J107").Value = 7188.5
$ws.Range("K107").Value = 6590.1055
$ws.Range("L107").Value = 7188.5
$ws.Range("M107").Value = -4670.1055
$ws.Range("N107").Value = -11028.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3943.6667
$ws.Range("I134").Value = 2899.6
$ws.Range("K134").Value = 8698.799999999999
$ws.Range("M134").Value = -6163.799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 20000
$ws.Range("I4").Value = 20000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 20000
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -19888
$ws.Range("N4").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 499
$ws.Range("J7").Value = 499
$ws.Range("L7").Value = 499
$ws.Range("N7").Value = -725

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1553.3823
$ws.Range("I31").Value = 1337.9688
$ws.Range("K31").Value = 1337.9688
$ws.Range("M31").Value = -1042.9688

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1553.3823
$ws.Range("I34").Value = 1337.9688
$ws.Range("K34").Value = 1337.9688
$ws.Range("M34").Value = -1135.9688

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 62503824
$ws.Range("J62").Value = 100003760
$ws.Range("L62").Value = 100003760
$ws.Range("N62").Value = -100005008

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 62503824
$ws.Range("J65").Value = 100003760
$ws.Range("L65").Value = 500018800
$ws.Range("N65").Value = -500025040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 9758975
$ws.Range("I99").Value = 1527320.4
$ws.Range("J99").Value = 28574186
$ws.Range("K99").Value = 1527320.4
$ws.Range("L99").Value = 28574186
$ws.Range("M99").Value = -1525822.4
$ws.Range("N99").Value = -28577182

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 886.4286
$ws.Range("I105").Value = 845.9
$ws.Range("J105").Value = 987.75
$ws.Range("K105").Value = 845.9
$ws.Range("L105").Value = 987.75
$ws.Range("M105").Value = 901.1
$ws.Range("N105").Value = -4481.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 9758975
$ws.Range("I126").Value = 1527320.4
$ws.Range("J126").Value = 28574186
$ws.Range("K126").Value = 4581961.199999999
$ws.Range("L126").Value = 85722558
$ws.Range("M126").Value = -4579491.199999999
$ws.Range("N126").Value = -85727498

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 865.7857
$ws.Range("I60").Value = 112.818184
$ws.Range("K60").Value = 338.454552
$ws.Range("M60").Value = -87.45455200000004

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2938
$ws.Range("J75").Value = 4605
$ws.Range("L75").Value = 13815
$ws.Range("N75").Value = -15811

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 2938
$ws.Range("J78").Value = 4605
$ws.Range("L78").Value = 41445
$ws.Range("N78").Value = -51429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 2624.9062
$ws.Range("J107").Value = 2591.318
$ws.Range("L107").Value = 7773.954000000001
$ws.Range("N107").Value = -11613.954

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 4166.6665
$ws.Range("J116").Value = 4166.6665
$ws.Range("L116").Value = 12499.9995
$ws.Range("N116").Value = -19383.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 4751.4546
$ws.Range("I140").Value = 3066.5715
$ws.Range("K140").Value = 9199.7145
$ws.Range("M140").Value = -4019.7145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5363.3335
$ws.Range("J80").Value = 7231.222
$ws.Range("L80").Value = 7231.222
$ws.Range("N80").Value = -9227.222

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 5363.3335
$ws.Range("J83").Value = 7231.222
$ws.Range("L83").Value = 36156.11
$ws.Range("N83").Value = -46140.11

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 6789.2354
$ws.Range("J113").Value = 8680.125
$ws.Range("L113").Value = 8680.125
$ws.Range("N113").Value = -13020.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3373.25
$ws.Range("I126").Value = 2997.2
$ws.Range("K126").Value = 8991.599999999999
$ws.Range("M126").Value = -6521.599999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 471.68
$ws.Range("I55").Value = 481.2857
$ws.Range("J55").Value = 421.25
$ws.Range("K55").Value = 481.2857
$ws.Range("L55").Value = 421.25
$ws.Range("M55").Value = -308.2857
$ws.Range("N55").Value = -767.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H128").Value = 49999
$ws.Range("J128").Value = 49999
$ws.Range("L128").Value = 49999
$ws.Range("N128").Value = -59959

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2832.7036
$ws.Range("J132").Value = 3351.4666
$ws.Range("L132").Value = 10054.3998
$ws.Range("N132").Value = -15114.3998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 18777.777
$ws.Range("J109").Value = 18777.777
$ws.Range("L109").Value = 18777.777
$ws.Range("N109").Value = -21551.777

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2042.5294
$ws.Range("J122").Value = 3999.25
$ws.Range("L122").Value = 11997.75
$ws.Range("N122").Value = -16897.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 34373.375
$ws.Range("J125").Value = 33172.22
$ws.Range("L125").Value = 33172.22
$ws.Range("N125").Value = -43012.22

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3652.6511
$ws.Range("I132").Value = 3292.7646
$ws.Range("K132").Value = 9878.293799999999
$ws.Range("M132").Value = -7348.293799999999
